$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '42.641.22'
$ws.Range('E2').Value = '  -0.49%  '
$ws.Range('D3').Value = '2.532.60'
$ws.Range('E3').Value = '  -1.43%  '
$ws.Range('D4').Value = '''1.00'
$ws.Range('E4').Value = '  +0.01%  '
$ws.Range('D5').Value = '''309.35'
$ws.Range('E5').Value = '  -1.53%  '
$ws.Range('D6').Value = '''100.14'
$ws.Range('E6').Value = '  +0.40%  '
$ws.Range('E7').Value = '  -1.40%  '
$ws.Range('D9').Value = '''0.522'
$ws.Range('E9').Value = '  -2.19%  '
$ws.Range('D10').Value = '''35.74'
$ws.Range('E10').Value = '  -0.83%  '
$ws.Range('D11').Value = '''0.0804'
$ws.Range('E11').Value = '  -0.74%  '
$ws.Range('D12').Value = '''7.36'
$ws.Range('E12').Value = '  -1.74%  '
$ws.Range('E13').Value = '  +0.66%  '
$ws.Range('D14').Value = '2.920.01'
$ws.Range('E14').Value = '  -1.59%  '
$ws.Range('D15').Value = '''15.32'
$ws.Range('E15').Value = '  -3.28%  '
$ws.Range('D16').Value = '2.474.89'
$ws.Range('E16').Value = '  -2.93%  '
$ws.Range('D17').Value = '''0.813'
$ws.Range('E17').Value = '  -4.10%  '
$ws.Range('D18').Value = '42.631.28'
$ws.Range('E18').Value = '  -0.64%  '
$ws.Range('D19').Value = '''6.72'
$ws.Range('E19').Value = '  -0.99%  '
$ws.Range('E20').Value = '  -1.55%  '
$ws.Range('D21').Value = '''12.27'
$ws.Range('E21').Value = '  -2.43%  '
$ws.Range('E22').Value = '  -0.25%  '
$ws.Range('D23').Value = '''243.23'
$ws.Range('E23').Value = '  -2.47%  '
$ws.Range('E24').Value = '  -3.29%  '
$ws.Range('E25').Value = '  -1.99%  '
$ws.Range('E26').Value = '  +0.61%  '
$ws.Range('D27').Value = '''25.60'
$ws.Range('E27').Value = '  -5.12%  '
$ws.Range('E28').Value = '  -1.14%  '
$ws.Range('D29').Value = '''10.17'
$ws.Range('E29').Value = '  -0.90%  '
$ws.Range('D30').Value = '''38.52'
$ws.Range('E30').Value = '  -5.13%  '
$ws.Range('D31').Value = '''158.40'
$ws.Range('E31').Value = '  +0.56%  '
$ws.Range('D32').Value = '''5.77'
$ws.Range('E32').Value = '  -0.58%  '
$ws.Range('E33').Value = '  +9.74%  '
$ws.Range('E34').Value = '  -1.52%  '
$ws.Range('D35').Value = '''0.0782'
$ws.Range('E35').Value = '  -2.32%  '
$ws.Range('D36').Value = '''18.34'
$ws.Range('E36').Value = '  -2.78%  '
$ws.Range('D37').Value = '''3.14'
$ws.Range('E37').Value = '  -7.40%  '
$ws.Range('D38').Value = '''1.97'
$ws.Range('E38').Value = '  -7.13%  '
$ws.Range('E39').Value = '  -1.46%  '
$ws.Range('E40').Value = '  -0.89%  '
$ws.Range('E41').Value = '  +3.72%  '
$ws.Range('D42').Value = '''22.61'
$ws.Range('E42').Value = '  -3.33%  '
$ws.Range('E43').Value = '  +0.13%  '
$ws.Range('E44').Value = '  -1.07%  '
$ws.Range('E45').Value = '  +0.98%  '
$ws.Range('D46').Value = '1.987.14'
$ws.Range('E46').Value = '  -1.08%  '
$ws.Range('E47').Value = '  -0.25%  '
$ws.Range('D48').Value = '2.776.86'
$ws.Range('E48').Value = '  -1.47%  '
$ws.Range('B49').Value = 'Algorand'
$ws.Range('C49').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D49').Value = '''0.189'
$ws.Range('E49').Value = '  -3.41%  '
$ws.Range('B50').Value = 'BitcoinSV'
$ws.Range('C50').Value = 'https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv'
$ws.Range('D50').Value = '''79.43'
$ws.Range('E50').Value = '  -2.72%  '
$ws.Range('D51').Value = '''72.11'
$ws.Range('E51').Value = '  -3.49%  '
